# Actualización automática 2025-07-24 17:25:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Cells.Item(11, 13).Value = 5.94          # M11: PORCELANATO sale for DDH S.A.S.
$wsGrupo.Cells.Item(24, 13).Value = "6 de 22"     # M24: count of advisors with PORCELANATO sales

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Cells.Item(11, 6).Value = 5.94         # F11: julio sale for DDH S.A.S.
$wsMensual.Cells.Item(24, 6).Value = 30727.2      # F24: julio total

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Cells.Item(16, 4).Value = 27190.22             # D16: PORCELANATO VENTA
$wsCumplimiento.Cells.Item(16, 5).Value = 11566.32             # E16: PORCELANATO POR CUMPLIR
$wsCumplimiento.Cells.Item(16, 6).Value = 0.7015646907592886   # F16: PORCELANATO CUMPLIMIENTO
$wsCumplimiento.Cells.Item(19, 4).Value = 30727.2              # D19: TOTAL VENTA
$wsCumplimiento.Cells.Item(19, 5).Value = 27495.80386304604    # E19: TOTAL POR CUMPLIR
$wsCumplimiento.Cells.Item(19, 6).Value = 0.5277501667945109   # F19: TOTAL CUMPLIMIENTO
